$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the paragraph that contains $anchorText (a unique substring)
# and replace that whole paragraph (including its end-of-paragraph mark) with
# the supplied OOXML for a single <w:p>...</w:p> (or <w:p/>) element.
# ---------------------------------------------------------------------------
function Replace-ParagraphXml($anchorText, $newParaXml) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng.Expand(4) | Out-Null
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
        '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1. "Git Log" paragraph: split the long run into several runs, with
#    w:proofErr (gramStart/gramEnd/spellStart/spellEnd) markers interleaved
#    -- same visible text, just re-run the proofing split.
# ---------------------------------------------------------------------------
$para15 = '<w:p w14:paraId="223B78F5" w14:textId="26B04EAA" w:rsidR="005B7310" w:rsidRDefault="00F91501" w:rsidP="00F91501">' + `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr>' + `
    '<w:r w:rsidRPr="00F91501"><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="single"/></w:rPr><w:t>Git Log</w:t></w:r>' + `
    '<w:r w:rsidRPr="00F91501"><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' + `
    '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">A git log </w:t></w:r>' + `
    '<w:r w:rsidR="007D7EDA"><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">shows lists of </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>commits(</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">SHA) that was made to repository, each commits associated with </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>message,commitSha,Author,date</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>.</w:t></w:r>' + `
    '</w:p>'
Replace-ParagraphXml "shows lists of" $para15

# ---------------------------------------------------------------------------
# 2. "`git log - -author=...`" paragraph: split run, proofErr gramStart/gramEnd
# ---------------------------------------------------------------------------
$para18 = '<w:p w14:paraId="36BB41BA" w14:textId="1A3D60F2" w:rsidR="007D7EDA" w:rsidRPr="007D7EDA" w:rsidRDefault="007D7EDA" w:rsidP="007D7EDA">' + `
    '<w:pPr><w:ind w:left="360"/><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">To check the author who committed changes into the repository - </w:t></w:r>' + `
    '<w:r w:rsidRPr="007D7EDA"><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>`git log - -author</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>=”Roni</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>”</w:t></w:r>' + `
    '<w:r w:rsidR="00445225"><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>.</w:t></w:r>' + `
    '</w:p>'
Replace-ParagraphXml "git log - -author" $para18

# ---------------------------------------------------------------------------
# 3. "`git log - -committer=...`" paragraph: split run, proofErr gramStart/gramEnd
# ---------------------------------------------------------------------------
$para19 = '<w:p w14:paraId="70BFB61C" w14:textId="3B6C17BB" w:rsidR="007D7EDA" w:rsidRDefault="007D7EDA" w:rsidP="007D7EDA">' + `
    '<w:pPr><w:ind w:left="360"/><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">To check the committer who committed into the repository - </w:t></w:r>' + `
    '<w:r w:rsidRPr="007D7EDA"><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>`git log - -committer</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>=”Roni</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>"`</w:t></w:r>' + `
    '<w:r w:rsidR="004758A5"><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>.</w:t></w:r>' + `
    '</w:p>'
Replace-ParagraphXml "git log - -committer" $para19

# ---------------------------------------------------------------------------
# 4. Append a new empty paragraph plus the "cherry pick" paragraph at the
#    very end of the document body (before sectPr).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRng = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newTail = '<w:p/>' + `
    '<w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="single"/></w:rPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">To cherry pick the commits from one branch to another, branch1 copy sha by using git log then switch to current working head branch and do - </w:t></w:r>' + `
    '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>`git cherry-pick 2acsacs -&gt; git push`</w:t></w:r>' + `
    '</w:p>'

$pkgTail = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:body>' + $newTail + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRng.InsertXML($pkgTail)

Write-Host "done"
